$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 17 / E17: length value corrected from 20 to 30
$ws.Range("E17").Value = "30"

# Row 17 / G17: note explaining the change, highlighted in yellow
$note = "2022-03-11 智偉修改`n長度原20新30"
$ws.Range("G17").Value = $note
$ws.Range("G17").WrapText = $true
$ws.Range("G17").HorizontalAlignment = -4131
$ws.Range("G17").VerticalAlignment = -4160
$ws.Range("G17").Interior.Color = 65535

# Only style the "智偉修改\n長度原20新30" portion of the note with the 細明體 font,
# leaving the leading date stamp in the default font (mirrors the source edit).
$dateLen = ("2022-03-11 ").Length
$chars = $ws.Range("G17").Characters($dateLen + 1, $note.Length - $dateLen)
$chars.Font.Name = "細明體"
$chars.Font.Size = 12

# Row 17 grows to fit the two-line note
$ws.Rows.Item(17).RowHeight = 32.4

# Leave the selection where the author left it after making the edit
[void]$ws.Range("C9").Select()
